# Fruta / hortaliza, semanal
# Insert a new weekly record at row 13, pushing the existing rows 13-16 down
# to rows 14-17 (dimension grows from A1:T16 to A1:T17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 13 (this shifts old rows
# 13-16 down to 14-17, copying formatting from the row above as Excel does).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new weekly data point.
$ws.Range("A13").Value2 = 1
$ws.Range("B13").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value2 = "Arica y Parinacota"
$ws.Range("D13").Value2 = 44784
$ws.Range("E13").Value2 = 15
$ws.Range("F13").Value2 = "Fruta"
$ws.Range("G13").Value2 = 100101
$ws.Range("H13").Value2 = "Berries"
$ws.Range("I13").Value2 = 100101007
$ws.Range("J13").Value2 = "Kiwi"
$ws.Range("K13").Value2 = "Hayward"
$ws.Range("L13").Value2 = "Primera"
$ws.Range("M13").Value2 = 300
$ws.Range("N13").Value2 = 19000
$ws.Range("O13").Value2 = 20000
$ws.Range("P13").Value2 = 19500
$ws.Range("Q13").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R13").Value2 = "Región de O'Higgins"
$ws.Range("S13").Value2 = 1083
$ws.Range("T13").Value2 = 18
